# Weekly update: insert two new price records (Apio, Terminal Hortofrutícola
# Agro Chillán) at the top of the data block (rows 278-279), pushing the
# existing rows down by two. This mirrors the "fruta / hortaliza, semanal"
# refresh: the newest week's observations are inserted right after the most
# recent existing row, and every older row slides down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the current row 278.
$ws.Rows.Item(278).Insert()
$ws.Rows.Item(278).Insert()

# New row 278: "Primera" quality record for the new date (45093).
$ws.Cells.Item(278, 1).Value = 7
$ws.Cells.Item(278, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(278, 3).Value = "Ñuble"
$ws.Cells.Item(278, 4).Value = 45093
$ws.Cells.Item(278, 5).Value = 16
$ws.Cells.Item(278, 6).Value = 100112017
$ws.Cells.Item(278, 7).Value = "Apio"
$ws.Cells.Item(278, 8).Value = "Americana (o)"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 80
$ws.Cells.Item(278, 11).Value = 7000
$ws.Cells.Item(278, 12).Value = 7000
$ws.Cells.Item(278, 13).Value = 7000
$ws.Cells.Item(278, 14).Value = "$/docena de matas"
$ws.Cells.Item(278, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(278, 16).Value = 1167
$ws.Cells.Item(278, 17).Value = 6
$ws.Cells.Item(278, 18).Value = "Hortaliza"

# New row 279: "Segunda" quality record for the same new date (45093).
$ws.Cells.Item(279, 1).Value = 7
$ws.Cells.Item(279, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(279, 3).Value = "Ñuble"
$ws.Cells.Item(279, 4).Value = 45093
$ws.Cells.Item(279, 5).Value = 16
$ws.Cells.Item(279, 6).Value = 100112017
$ws.Cells.Item(279, 7).Value = "Apio"
$ws.Cells.Item(279, 8).Value = "Americana (o)"
$ws.Cells.Item(279, 9).Value = "Segunda"
$ws.Cells.Item(279, 10).Value = 150
$ws.Cells.Item(279, 11).Value = 6000
$ws.Cells.Item(279, 12).Value = 6000
$ws.Cells.Item(279, 13).Value = 6000
$ws.Cells.Item(279, 14).Value = "$/docena de matas"
$ws.Cells.Item(279, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(279, 16).Value = 1000
$ws.Cells.Item(279, 17).Value = 6
$ws.Cells.Item(279, 18).Value = "Hortaliza"
